$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StringKey")

# Row 11 - new HUD interaction string for missions
$ws.Cells.Item(11, 1).Value = "sys.hud.interaction.mission"
$ws.Cells.Item(11, 2).Value = "F를 눌러 미션 시작: {0}"
$ws.Cells.Item(11, 3).Value = "en9"

# Column A (keys) for rows 12-21
$keys = @(
    "sys.mission.name.1",
    "sys.mission.name.2",
    "sys.mission.name.3",
    "sys.mission.name.4",
    "sys.mission.name.5",
    "sys.mission.content.1",
    "sys.mission.content.2",
    "sys.mission.content.3",
    "sys.mission.content.4",
    "sys.mission.content.5"
)
for ($i = 0; $i -lt $keys.Length; $i++) {
    $ws.Cells.Item(12 + $i, 1).Value = $keys[$i]
}

# Column B (Korean values) for rows 12-21
$krValues = @(
    "미션 1 이름입니다~~",
    "미션 2 이름입니다~~",
    "미션 3 이름입니다~~",
    "미션 4 이름입니다~~",
    "미션 5 이름입니다~~",
    "미션 1 내용입니다~~ 임시로 넣어놓은 긴 내용입니다.",
    "미션 2 내용입니다~~ 임시로 넣어놓은 긴 내용입니다.",
    "미션 3 내용입니다~~ 임시로 넣어놓은 긴 내용입니다.",
    "미션 4 내용입니다~~ 임시로 넣어놓은 긴 내용입니다.",
    "미션 5 내용입니다~~ 임시로 넣어놓은 긴 내용입니다."
)
for ($i = 0; $i -lt $krValues.Length; $i++) {
    $ws.Cells.Item(12 + $i, 2).Value = $krValues[$i]
}

# Column C (English placeholder values) for rows 12-21
for ($i = 10; $i -le 19; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = "en$i"
}
